# Applies the OOXML-level numeric restatement of rows 2-9 (company_list sheet)
# described by the commit "error solve ifrs list": every quarterly-period row
# is rescaled from raw KRW amounts to a smaller reporting unit, a handful of
# now-redundant sub-total cells are dropped, and the last two periods (rows 8-9)
# are removed entirely (kept only as row labels).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = 3062
$ws.Range("E2").Value = 568
$ws.Range("F2").Value = 568
$ws.Range("G2").Value = 600
$ws.Range("H2").Value = 466
$ws.Range("I2").Value = 440
$ws.Range("J2").Value = 26
$ws.Range("K2").Value = 6412
$ws.Range("L2").Value = 680
$ws.Range("M2").Value = 5732
$ws.Range("N2").Value = 5410
$ws.Range("O2").Value = 322
$ws.Range("P2").Value = 539
$ws.Range("Q2").Value = 1017
$ws.Range("R2").Value = -953
$ws.Range("S2").Value = -187
$ws.Range("T2").Value = 315
$ws.Range("U2").Value = 702
$ws.Range("V2").Value = 120
$ws.Range("W2").Value = 18.55
$ws.Range("X2").Value = 15.22
$ws.Range("Y2").Value = 8.42
$ws.Range("Z2").Value = 7.44
$ws.Range("AA2").Value = 11.85
$ws.Range("AB2").Value = 965.98
$ws.Range("AC2").Value = 408
$ws.Range("AD2").Value = 10.87
$ws.Range("AE2").Value = 5016
$ws.Range("AF2").Value = 0.88
$ws.Range("AG2").Value = 40
$ws.Range("AH2").Value = 0.9
$ws.Range("AI2").Value = 9.800000000000001
$ws.Range("AJ2").Value = 107892760

# --- Row 3 ---
$ws.Range("D3").Value = 2912
$ws.Range("E3").Value = 454
$ws.Range("F3").Value = 454
$ws.Range("G3").Value = 500
$ws.Range("H3").Value = 382
$ws.Range("I3").Value = 370
$ws.Range("J3").Value = 12
$ws.Range("K3").Value = 6522
$ws.Range("L3").Value = 556
$ws.Range("M3").Value = 5966
$ws.Range("N3").Value = 5774
$ws.Range("O3").Value = 191
$ws.Range("P3").Value = 539
$ws.Range("Q3").Value = 813
$ws.Range("R3").Value = -452
$ws.Range("S3").Value = -261
$ws.Range("T3").Value = 345
$ws.Range("U3").Value = 468
$ws.Range("W3").Value = 15.58
$ws.Range("X3").Value = 13.12
$ws.Range("Y3").Value = 6.61
$ws.Range("Z3").Value = 5.91
$ws.Range("AA3").Value = 9.32
$ws.Range("AB3").Value = 1025.19
$ws.Range("AC3").Value = 343
$ws.Range("AD3").Value = 10.15
$ws.Range("AE3").Value = 5355
$ws.Range("AF3").Value = 0.65
$ws.Range("AG3").Value = 40
$ws.Range("AH3").Value = 1.15
$ws.Range("AI3").Value = 11.66
$ws.Range("AJ3").Value = 107892760
$ws.Range("V3").ClearContents()

# --- Row 4 ---
$ws.Range("D4").Value = 2921
$ws.Range("E4").Value = 492
$ws.Range("F4").Value = 492
$ws.Range("G4").Value = 525
$ws.Range("H4").Value = 456
$ws.Range("I4").Value = 445
$ws.Range("J4").Value = 11
$ws.Range("K4").Value = 6911
$ws.Range("L4").Value = 627
$ws.Range("M4").Value = 6284
$ws.Range("N4").Value = 6284
$ws.Range("P4").Value = 564
$ws.Range("Q4").Value = 965
$ws.Range("R4").Value = -530
$ws.Range("S4").Value = -139
$ws.Range("T4").Value = 203
$ws.Range("U4").Value = 761
$ws.Range("W4").Value = 16.86
$ws.Range("X4").Value = 15.61
$ws.Range("Y4").Value = 7.37
$ws.Range("Z4").Value = 6.79
$ws.Range("AA4").Value = 9.98
$ws.Range("AB4").Value = 1081.9
$ws.Range("AC4").Value = 412
$ws.Range("AD4").Value = 10.04
$ws.Range("AE4").Value = 5569
$ws.Range("AF4").Value = 0.74
$ws.Range("AG4").Value = 40
$ws.Range("AH4").Value = 0.97
$ws.Range("AI4").Value = 9.91
$ws.Range("AJ4").Value = 112876596
$ws.Range("O4").ClearContents()
$ws.Range("V4").ClearContents()

# --- Row 5 ---
$ws.Range("D5").Value = 2902
$ws.Range("E5").Value = 493
$ws.Range("F5").Value = 493
$ws.Range("G5").Value = 563
$ws.Range("H5").Value = 435
$ws.Range("I5").Value = 435
$ws.Range("K5").Value = 7350
$ws.Range("L5").Value = 658
$ws.Range("M5").Value = 6691
$ws.Range("N5").Value = 6691
$ws.Range("P5").Value = 564
$ws.Range("Q5").Value = 872
$ws.Range("R5").Value = -1143
$ws.Range("S5").Value = -44
$ws.Range("T5").Value = 240
$ws.Range("U5").Value = 632
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 16.97
$ws.Range("X5").Value = 15
$ws.Range("Y5").Value = 6.71
$ws.Range("Z5").Value = 6.11
$ws.Range("AA5").Value = 9.84
$ws.Range("AB5").Value = 1153.75
$ws.Range("AC5").Value = 386
$ws.Range("AD5").Value = 10.24
$ws.Range("AE5").Value = 6075
$ws.Range("AF5").Value = 0.65
$ws.Range("AG5").Value = 50
$ws.Range("AH5").Value = 1.27
$ws.Range("AI5").Value = 12.65
$ws.Range("AJ5").Value = 112876596
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# --- Row 6 ---
$ws.Range("D6").Value = 2897
$ws.Range("E6").Value = 467
$ws.Range("F6").Value = 467
$ws.Range("G6").Value = 529
$ws.Range("H6").Value = 399
$ws.Range("I6").Value = 399
$ws.Range("K6").Value = 7805
$ws.Range("L6").Value = 664
$ws.Range("M6").Value = 7140
$ws.Range("N6").Value = 7140
$ws.Range("P6").Value = 564
$ws.Range("Q6").Value = 632
$ws.Range("R6").Value = -795
$ws.Range("S6").Value = -55
$ws.Range("T6").Value = 140
$ws.Range("U6").Value = 493
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 16.13
$ws.Range("X6").Value = 13.78
$ws.Range("Y6").Value = 5.77
$ws.Range("Z6").Value = 5.27
$ws.Range("AA6").Value = 9.300000000000001
$ws.Range("AB6").Value = 1233.26
$ws.Range("AC6").Value = 354
$ws.Range("AD6").Value = 11.48
$ws.Range("AE6").Value = 6482
$ws.Range("AF6").Value = 0.63
$ws.Range("AG6").Value = 60
$ws.Range("AH6").Value = 1.48
$ws.Range("AI6").Value = 16.55
$ws.Range("AJ6").Value = 112876596

# --- Row 7 ---
$ws.Range("D7").Value = 2970
$ws.Range("E7").Value = 430
$ws.Range("I7").Value = 390
$ws.Range("W7").Value = 14.48
$ws.Range("AC7").Value = 346
$ws.Range("AD7").Value = 11
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# --- Row 8 ---
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# --- Row 9 ---
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()

Write-Output "Applied restated ifrs figures to rows 2-9"
